$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (39 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4363390
$ws.Range("I86").Value = 7167583
$ws.Range("J86").Value = 1311.5555
$ws.Range("K86").Value = 7167583
$ws.Range("L86").Value = 1311.5555
$ws.Range("M86").Value = -7166460
$ws.Range("N86").Value = -3557.5555
$ws.Range("H89").Value = 4363390
$ws.Range("I89").Value = 7167583
$ws.Range("J89").Value = 1311.5555
$ws.Range("K89").Value = 35837915
$ws.Range("L89").Value = 6557.7775
$ws.Range("M89").Value = -35832299
$ws.Range("N89").Value = -17789.7775
$ws.Range("H92").Value = 8773183
$ws.Range("I92").Value = 11494998
$ws.Range("J92").Value = 2887.111
$ws.Range("K92").Value = 11494998
$ws.Range("L92").Value = 2887.111
$ws.Range("M92").Value = -11493750
$ws.Range("N92").Value = -5383.111
$ws.Range("H107").Value = 1282.2307
$ws.Range("I107").Value = 1603.9375
$ws.Range("J107").Value = 767.5
$ws.Range("K107").Value = 1603.9375
$ws.Range("L107").Value = 767.5
$ws.Range("M107").Value = 316.0625
$ws.Range("N107").Value = -4607.5
$ws.Range("H112").Value = 12293.479
$ws.Range("J112").Value = 12293.479
$ws.Range("L112").Value = 36880.437
$ws.Range("N112").Value = -39096.437
$ws.Range("H137").Value = 1567.65
$ws.Range("I137").Value = 1490.2222
$ws.Range("J137").Value = 1631
$ws.Range("K137").Value = 4470.6666
$ws.Range("L137").Value = 4893
$ws.Range("M137").Value = -1920.6666
$ws.Range("N137").Value = -9993

# --- Sheet: ARM (7 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2100
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 200
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = 296
$ws.Range("N97").Value = -4992

# --- Sheet: BSM (4 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3841.5789
$ws.Range("I134").Value = 3811.875
$ws.Range("K134").Value = 11435.625
$ws.Range("M134").Value = -8900.625

# --- Sheet: CRP (8 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4601.222
$ws.Range("I31").Value = 1166.4517
$ws.Range("K31").Value = 1166.4517
$ws.Range("M31").Value = -871.4517000000001
$ws.Range("H34").Value = 4601.222
$ws.Range("I34").Value = 1166.4517
$ws.Range("K34").Value = 1166.4517
$ws.Range("M34").Value = -964.4517000000001

# --- Sheet: CUL (74 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 479.35715
$ws.Range("I15").Value = 444.4
$ws.Range("J15").Value = 486.9565
$ws.Range("K15").Value = 1333.2
$ws.Range("L15").Value = 1460.8695
$ws.Range("M15").Value = -1193.2
$ws.Range("N15").Value = -1740.8695
$ws.Range("H20").Value = 947.8261
$ws.Range("H21").Value = 2875.125
$ws.Range("J21").Value = 3271.4285
$ws.Range("L21").Value = 9814.2855
$ws.Range("N21").Value = -10160.2855
$ws.Range("H22").Value = 842.2692
$ws.Range("I22").Value = 544.3333
$ws.Range("K22").Value = 1632.9999
$ws.Range("M22").Value = -1463.9999
$ws.Range("H26").Value = 6124.3145
$ws.Range("I26").Value = 81.625
$ws.Range("J26").Value = 7914.7407
$ws.Range("K26").Value = 244.875
$ws.Range("L26").Value = 23744.2221
$ws.Range("M26").Value = 43.125
$ws.Range("N26").Value = -24320.2221
$ws.Range("H27").Value = 842.2692
$ws.Range("I27").Value = 544.3333
$ws.Range("K27").Value = 1632.9999
$ws.Range("M27").Value = -1530.9999
$ws.Range("H32").Value = 47623336
$ws.Range("J32").Value = 47623336
$ws.Range("L32").Value = 142870008
$ws.Range("N32").Value = -142870574
$ws.Range("H34").Value = 16129567
$ws.Range("I34").Value = 146.25
$ws.Range("J34").Value = 18519112
$ws.Range("K34").Value = 438.75
$ws.Range("L34").Value = 55557336
$ws.Range("M34").Value = -354.75
$ws.Range("N34").Value = -55557504
$ws.Range("H39").Value = 3150
$ws.Range("J39").Value = 3150
$ws.Range("L39").Value = 9450
$ws.Range("N39").Value = -10038
$ws.Range("H44").Value = 39931.223
$ws.Range("I44").Value = 71320
$ws.Range("J44").Value = 695.25
$ws.Range("K44").Value = 213960
$ws.Range("L44").Value = 2085.75
$ws.Range("M44").Value = -213562
$ws.Range("N44").Value = -2881.75
$ws.Range("H46").Value = 1515.7142
$ws.Range("I46").Value = 333.33334
$ws.Range("J46").Value = 1988.6666
$ws.Range("K46").Value = 1000.00002
$ws.Range("L46").Value = 5965.9998
$ws.Range("M46").Value = -909.0000200000001
$ws.Range("N46").Value = -6147.9998
$ws.Range("H57").Value = 1595
$ws.Range("I57").Value = 433.33334
$ws.Range("J57").Value = 1800
$ws.Range("K57").Value = 1300.00002
$ws.Range("L57").Value = 5400
$ws.Range("M57").Value = -741.0000199999999
$ws.Range("N57").Value = -6518
$ws.Range("H58").Value = 901.6667
$ws.Range("I58").Value = 682
$ws.Range("K58").Value = 2046
$ws.Range("M58").Value = -1918
$ws.Range("H131").Value = 4109.9736
$ws.Range("I131").Value = 522.8570999999999
$ws.Range("J131").Value = 4919.968
$ws.Range("K131").Value = 1568.5713
$ws.Range("L131").Value = 14759.904
$ws.Range("M131").Value = 3471.4287
$ws.Range("N131").Value = -24839.904

# --- Sheet: GSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2107.8572
$ws.Range("I97").Value = 2031.25
$ws.Range("J97").Value = 2210
$ws.Range("K97").Value = 2031.25
$ws.Range("L97").Value = 2210
$ws.Range("M97").Value = -1535.25
$ws.Range("N97").Value = -3202
$ws.Range("H139").Value = 63946.92
$ws.Range("J139").Value = 63946.92
$ws.Range("L139").Value = 63946.92
$ws.Range("N139").Value = -74226.92

# --- Sheet: LTW (21 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2390.182
$ws.Range("I68").Value = 2198.4
$ws.Range("J68").Value = 2550
$ws.Range("K68").Value = 2198.4
$ws.Range("L68").Value = 2550
$ws.Range("M68").Value = -1449.4
$ws.Range("N68").Value = -4048
$ws.Range("H71").Value = 2390.182
$ws.Range("I71").Value = 2198.4
$ws.Range("J71").Value = 2550
$ws.Range("K71").Value = 10992
$ws.Range("L71").Value = 12750
$ws.Range("M71").Value = -7248
$ws.Range("N71").Value = -20238
$ws.Range("H93").Value = 12689.223
$ws.Range("I93").Value = 15314.714
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 15314.714
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = -14066.714
$ws.Range("N93").Value = -5996

Write-Host "Applied all Anima_Profits updates"